$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 820, shifting existing data (rows 820-870) down to 823-873
$ws.Range("A820:A822").EntireRow.Insert()

# Row 820
$ws.Cells.Item(820, 1).Value = 3
$ws.Cells.Item(820, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(820, 3).Value = "Coquimbo"
$ws.Cells.Item(820, 4).Value = 44706
$ws.Cells.Item(820, 5).Value = 5
$ws.Cells.Item(820, 6).Value = "Fruta"
$ws.Cells.Item(820, 7).Value = 100108
$ws.Cells.Item(820, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(820, 9).Value = 100108006
$ws.Cells.Item(820, 10).Value = "Plátano"
$ws.Cells.Item(820, 11).Value = "Sin especificar"
$ws.Cells.Item(820, 12).Value = "Maduro"
$ws.Cells.Item(820, 13).Value = 160
$ws.Cells.Item(820, 14).Value = 11000
$ws.Cells.Item(820, 15).Value = 11000
$ws.Cells.Item(820, 16).Value = 11000
$ws.Cells.Item(820, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(820, 18).Value = "Ecuador"
$ws.Cells.Item(820, 19).Value = 550
$ws.Cells.Item(820, 20).Value = 20

# Row 821
$ws.Cells.Item(821, 1).Value = 3
$ws.Cells.Item(821, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(821, 3).Value = "Coquimbo"
$ws.Cells.Item(821, 4).Value = 44706
$ws.Cells.Item(821, 5).Value = 5
$ws.Cells.Item(821, 6).Value = "Fruta"
$ws.Cells.Item(821, 7).Value = 100108
$ws.Cells.Item(821, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(821, 9).Value = 100108006
$ws.Cells.Item(821, 10).Value = "Plátano"
$ws.Cells.Item(821, 11).Value = "Sin especificar"
$ws.Cells.Item(821, 12).Value = "Pintón"
$ws.Cells.Item(821, 13).Value = 240
$ws.Cells.Item(821, 14).Value = 12000
$ws.Cells.Item(821, 15).Value = 12000
$ws.Cells.Item(821, 16).Value = 12000
$ws.Cells.Item(821, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(821, 18).Value = "Ecuador"
$ws.Cells.Item(821, 19).Value = 600
$ws.Cells.Item(821, 20).Value = 20

# Row 822
$ws.Cells.Item(822, 1).Value = 3
$ws.Cells.Item(822, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(822, 3).Value = "Coquimbo"
$ws.Cells.Item(822, 4).Value = 44706
$ws.Cells.Item(822, 5).Value = 5
$ws.Cells.Item(822, 6).Value = "Fruta"
$ws.Cells.Item(822, 7).Value = 100108
$ws.Cells.Item(822, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(822, 9).Value = 100108006
$ws.Cells.Item(822, 10).Value = "Plátano"
$ws.Cells.Item(822, 11).Value = "Sin especificar"
$ws.Cells.Item(822, 12).Value = "Primera Pintón"
$ws.Cells.Item(822, 13).Value = 280
$ws.Cells.Item(822, 14).Value = 13000
$ws.Cells.Item(822, 15).Value = 13000
$ws.Cells.Item(822, 16).Value = 13000
$ws.Cells.Item(822, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(822, 18).Value = "Ecuador"
$ws.Cells.Item(822, 19).Value = 650
$ws.Cells.Item(822, 20).Value = 20
